$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 31130
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 31130
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 31130
$ws.Range("N3").Value = -31358

$ws.Range("H11").Value = 13.705882
$ws.Range("I11").Value = 13.705882
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 13.705882
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 126.294118

$ws.Range("H33").Value = 88.875
$ws.Range("I33").Value = 88.875
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 88.875
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 140.125

$ws.Range("H70").Value = 9112.857
$ws.Range("I70").Value = 8433.333000000001
$ws.Range("J70").Value = 9622.5
$ws.Range("K70").Value = 25299.999
$ws.Range("L70").Value = 28867.5
$ws.Range("M70").Value = -25029.999
$ws.Range("N70").Value = -29407.5

$ws.Range("H73").Value = 9112.857
$ws.Range("I73").Value = 8433.333000000001
$ws.Range("J73").Value = 9622.5
$ws.Range("K73").Value = 25299.999
$ws.Range("L73").Value = 28867.5
$ws.Range("M73").Value = -24363.999
$ws.Range("N73").Value = -30739.5

$ws.Range("H95").Value = 23999
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 23999
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 23999
$ws.Range("N95").Value = -29491

$ws.Range("H102").Value = 31130
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 31130
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 31130
$ws.Range("N102").Value = -37620

$ws.Range("H131").Value = 2658.3333
$ws.Range("I131").Value = 2658.3333
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 7974.999899999999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -2934.999899999999

$ws.Range("H132").Value = 11940.211
$ws.Range("I132").Value = 11940.211
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 35820.633
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -33290.633

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 104500
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 104500
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 104500
$ws.Range("N140").Value = -114860

$ws.Range("H141").Value = 4733
$ws.Range("I141").Value = 4733
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 14199
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -9019

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6251168
$ws.Range("I32").Value = 1246.1333
$ws.Range("J32").Value = 100000000
$ws.Range("K32").Value = 1246.1333
$ws.Range("L32").Value = 100000000
$ws.Range("M32").Value = -959.1333

$ws.Range("H61").Value = 2476.375
$ws.Range("I61").Value = 1269.9375
$ws.Range("J61").Value = 4889.25
$ws.Range("K61").Value = 1269.9375
$ws.Range("L61").Value = 4889.25
$ws.Range("M61").Value = -1057.9375
$ws.Range("N61").Value = -5313.25

$ws.Range("H122").Value = 3042
$ws.Range("I122").Value = 1580.5
$ws.Range("J122").Value = 4503.5
$ws.Range("K122").Value = 4741.5
$ws.Range("L122").Value = 13510.5
$ws.Range("M122").Value = -2291.5
$ws.Range("N122").Value = -18410.5

$ws.Range("H132").Value = 1288.8
$ws.Range("I132").Value = 1239
$ws.Range("J132").Value = 1612.5
$ws.Range("K132").Value = 3717
$ws.Range("L132").Value = 4837.5
$ws.Range("M132").Value = -1187
$ws.Range("N132").Value = -9897.5

$ws.Range("H136").Value = 2476.375
$ws.Range("I136").Value = 1269.9375
$ws.Range("J136").Value = 4889.25
$ws.Range("K136").Value = 3809.8125
$ws.Range("L136").Value = 14667.75
$ws.Range("M136").Value = -1259.8125
$ws.Range("N136").Value = -19767.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5911.4614
$ws.Range("I86").Value = 2824.8333
$ws.Range("J86").Value = 8557.143
$ws.Range("K86").Value = 2824.8333
$ws.Range("L86").Value = 8557.143
$ws.Range("M86").Value = -1701.8333
$ws.Range("N86").Value = -10803.143

$ws.Range("H89").Value = 5911.4614
$ws.Range("I89").Value = 2824.8333
$ws.Range("J89").Value = 8557.143
$ws.Range("K89").Value = 14124.1665
$ws.Range("L89").Value = 42785.715
$ws.Range("M89").Value = -8508.166499999999
$ws.Range("N89").Value = -54017.715

$ws.Range("H99").Value = 83334530
$ws.Range("I99").Value = 100001220
$ws.Range("J99").Value = 1108.5
$ws.Range("K99").Value = 100001220
$ws.Range("L99").Value = 1108.5
$ws.Range("M99").Value = -99999722
$ws.Range("N99").Value = -4104.5

$ws.Range("H103").Value = 16000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 16000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 16000
$ws.Range("N103").Value = -18344

$ws.Range("H134").Value = 3711.739
$ws.Range("I134").Value = 1082.8948
$ws.Range("J134").Value = 16198.75
$ws.Range("K134").Value = 3248.6844
$ws.Range("L134").Value = 48596.25
$ws.Range("M134").Value = -713.6844000000001
$ws.Range("N134").Value = -53666.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 747
$ws.Range("I5").Value = 195.4
$ws.Range("J5").Value = 1666.3334
$ws.Range("K5").Value = 195.4
$ws.Range("L5").Value = 1666.3334
$ws.Range("M5").Value = -83.40000000000001
$ws.Range("N5").Value = -1890.3334

$ws.Range("H19").Value = 356.125
$ws.Range("I19").Value = 299.8
$ws.Range("J19").Value = 450
$ws.Range("K19").Value = 299.8
$ws.Range("L19").Value = 450
$ws.Range("M19").Value = -129.8
$ws.Range("N19").Value = -790

$ws.Range("H24").Value = 356.125
$ws.Range("I24").Value = 299.8
$ws.Range("J24").Value = 450
$ws.Range("K24").Value = 299.8
$ws.Range("L24").Value = 450
$ws.Range("M24").Value = -129.8
$ws.Range("N24").Value = -790

$ws.Range("H31").Value = 4874.189
$ws.Range("I31").Value = 2705.577
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 2705.577
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -2410.577

$ws.Range("H34").Value = 4874.189
$ws.Range("I34").Value = 2705.577
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 2705.577
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -2503.577

$ws.Range("H132").Value = 1431.4546
$ws.Range("I132").Value = 1431.4546
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4294.3638
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1764.3638

$ws.Range("H134").Value = 3268.0908
$ws.Range("I134").Value = 2523.0527
$ws.Range("J134").Value = 7986.6665
$ws.Range("K134").Value = 7569.158100000001
$ws.Range("L134").Value = 23959.9995
$ws.Range("M134").Value = -5034.158100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 908.8333
$ws.Range("I34").Value = 249
$ws.Range("J34").Value = 947.64703
$ws.Range("K34").Value = 747
$ws.Range("L34").Value = 2842.94109
$ws.Range("M34").Value = -663
$ws.Range("N34").Value = -3010.94109

$ws.Range("H39").Value = 5265.1665
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 5265.1665
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 15795.4995
$ws.Range("N39").Value = -16383.4995

$ws.Range("H55").Value = 1898.1666
$ws.Range("I55").Value = 425.125
$ws.Range("J55").Value = 3076.6
$ws.Range("K55").Value = 1275.375
$ws.Range("L55").Value = 9229.799999999999
$ws.Range("M55").Value = -1098.375
$ws.Range("N55").Value = -9583.799999999999

$ws.Range("H113").Value = 1523.55
$ws.Range("I113").Value = 1014
$ws.Range("J113").Value = 1797.9231
$ws.Range("K113").Value = 3042
$ws.Range("L113").Value = 5393.7693
$ws.Range("M113").Value = -872
$ws.Range("N113").Value = -9733.7693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1336

$ws.Range("H36").Value = 3005.6667
$ws.Range("I36").Value = 1758.5
$ws.Range("J36").Value = 5500
$ws.Range("K36").Value = 1758.5
$ws.Range("L36").Value = 5500
$ws.Range("M36").Value = -1273.5

$ws.Range("H132").Value = 47840.453
$ws.Range("I132").Value = 54946.105
$ws.Range("J132").Value = 2838
$ws.Range("K132").Value = 164838.315
$ws.Range("L132").Value = 8514
$ws.Range("M132").Value = -162308.315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5770.294
$ws.Range("I40").Value = 5630.9375
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 5630.9375
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -5494.9375
$ws.Range("N40").Value = -8272

$ws.Range("H82").Value = 3316.4614
$ws.Range("I82").Value = 1097.6
$ws.Range("J82").Value = 4703.25
$ws.Range("K82").Value = 1097.6
$ws.Range("L82").Value = 4703.25
$ws.Range("M82").Value = -736.5999999999999

$ws.Range("H85").Value = 3316.4614
$ws.Range("I85").Value = 1097.6
$ws.Range("J85").Value = 4703.25
$ws.Range("K85").Value = 1097.6
$ws.Range("L85").Value = 4703.25
$ws.Range("M85").Value = 150.4000000000001

$ws.Range("H122").Value = 3025.2856

$ws.Range("H132").Value = 9283.154
$ws.Range("I132").Value = 9425.546
$ws.Range("J132").Value = 8500
$ws.Range("K132").Value = 28276.638
$ws.Range("L132").Value = 25500
$ws.Range("M132").Value = -25746.638

$ws.Range("H136").Value = 967.3333
$ws.Range("I136").Value = 967.3333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2901.9999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -351.9998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 500
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 500
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -385
$ws.Range("N6").ClearContents()

$ws.Range("H122").Value = 5393.125
$ws.Range("I122").Value = 4211.25
$ws.Range("J122").Value = 6575
$ws.Range("K122").Value = 12633.75
$ws.Range("L122").Value = 19725
$ws.Range("M122").Value = -10183.75

$ws.Range("H132").Value = 947.6
$ws.Range("I132").Value = 908.1429000000001
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 2724.4287
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -194.4287000000004
$ws.Range("N132").Value = -9560

$ws.Range("H136").Value = 3029.15
$ws.Range("I136").Value = 2065.2222
$ws.Range("J136").Value = 3817.818
$ws.Range("K136").Value = 6195.6666
$ws.Range("L136").Value = 11453.454
$ws.Range("M136").Value = -3645.6666
$ws.Range("N136").Value = -16553.454
